$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 114, shifting the existing rows 114-115 down to 115-116
$ws.Rows.Item(114).Insert()

# Fill the new row 114 with the updated record
$ws.Range("A114").Value = 10
$ws.Range("B114").Value = "Vega Modelo de Temuco"
$ws.Range("C114").Value = "La Araucanía"
$ws.Range("D114").Value = 44448
$ws.Range("D114").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E114").Value = 9
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100104
$ws.Range("H114").Value = "Frutos de pepita"
$ws.Range("I114").Value = 100104003
$ws.Range("J114").Value = "Membrillo"
$ws.Range("K114").Value = "Champion"
$ws.Range("L114").Value = "Especial"
$ws.Range("M114").Value = 50
$ws.Range("N114").Value = 18000
$ws.Range("O114").Value = 19000
$ws.Range("P114").Value = 18400
$ws.Range("Q114").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R114").Value = "Región de O'Higgins"
$ws.Range("S114").Value = 1022
$ws.Range("T114").Value = 18
